$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers: lowercase "user" / "assistant"
$ws.Range("A1").Value = "user"
$ws.Range("B1").Value = "assistant"

# Row 2: Pablo Vegetti entry - fix "P." to "Pablo" in description
$ws.Range("A2").Value = "Pablo Vegetti"
$ws.Range("B2").Value = "Pablo Vegetti is a traditional target man, excelling in aerial duels and finishing chances in the box. He would be most effective in a system that relies on crossing and provides frequent service into the penalty area.
Vegetti's main strengths are: Aerial Ability, which indicates exceptional strength in aerial duels, making him a key target in the box< and Finishing and Opportunism, which suggets his finalization is strong and his opportunism are excellent, highlighting his clinical nature in scoring chances.
Vegetti's main weaknesses are: Playmaking, indicating limited contribution to creating chances for teammates; Defensive Work Rate, since his recomposition (defensive tracking back) is weak; and Dribbling: below-average performance in dribbling, reflecting limited ability to beat defenders individually.
Teams with a creative midfield and wingers delivering quality crosses would maximize his strengths, while his limited playmaking and defensive contributions would need to be offset by more dynamic teammates."

# Row 3: Gabriel Barbosa entry unchanged in content
$ws.Range("A3").Value = "Gabriel Barbosa"

# Update the view: select A2 (clears stale topLeftCell/selection state)
$ws.Range("A2").Select()
